# Modularized shopping cart class
# Update the credentials test-data row on Sheet1: replace the email/password
# pair in B2/C2, turn B2 into a mailto hyperlink (adding the built-in
# "Hyperlink" cell style along the way), and leave the selection on Sheet1
# at C2 (making Sheet1 the active/selected tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New credentials
$ws1.Range("B2").Value = "audi.love25@gmail.com"
$ws1.Range("C2").Value = "Mitsubishi7!"

# Turn the email cell into a real hyperlink (applies the Hyperlink style too)
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:audi.love25@gmail.com")

# Move the selection / active sheet to Sheet1!C2
$ws1.Activate()
$ws1.Range("C2").Select()
